$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 3 de Julio de 2020 a las 22:23"

# Update country data rows (values per official diff)
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 2875072
$ws.Range("C4").Value = 39388
$ws.Range("D4").Value = 1201749
$ws.Range("E4").Value = 1541389
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 449
$ws.Range("H4").Value = 131934

$ws.Range("A7").Value = "India"
$ws.Range("B7").Value = 649889
$ws.Range("C7").Value = 22721
$ws.Range("D7").Value = 394319
$ws.Range("E7").Value = 236901
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 444
$ws.Range("H7").Value = 18669

$ws.Range("A18").Value = "Alemania"
$ws.Range("B18").Value = 197000
$ws.Range("C18").Value = 283
$ws.Range("D18").Value = 181000
$ws.Range("E18").Value = 6927
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = 9073

$ws.Range("A26").Value = "Egipto"
$ws.Range("B26").Value = 72711
$ws.Range("C26").Value = 1412
$ws.Range("D26").Value = 19690
$ws.Range("E26").Value = 49820
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = 81
$ws.Range("H26").Value = 3201

$ws.Range("A35").Value = "Emiratos Arabes Unidos"
$ws.Range("B35").Value = 50141
$ws.Range("C35").Value = 672
$ws.Range("D35").Value = 39153
$ws.Range("E35").Value = 10670
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 1
$ws.Range("H35").Value = 318

$ws.Range("A41").Value = "Portugal"
$ws.Range("B41").Value = 43156
$ws.Range("C41").Value = 374
$ws.Range("D41").Value = 28424
$ws.Range("E41").Value = 13134
$ws.Range("F41").Value = 0
$ws.Range("G41").Value = 11
$ws.Range("H41").Value = 1598

$ws.Range("A70").Value = "Costa de Marfil"
$ws.Range("B70").Value = 10244
$ws.Range("C70").Value = 252
$ws.Range("D70").Value = 4726
$ws.Range("E70").Value = 5448
$ws.Range("F70").Value = 0
$ws.Range("G70").Value = 2
$ws.Range("H70").Value = 70

$ws.Range("A72").Value = "Uzbekistan"
$ws.Range("B72").Value = 9396
$ws.Range("C72").Value = 318
$ws.Range("D72").Value = 6251
$ws.Range("E72").Value = 3116
$ws.Range("F72").Value = 0
$ws.Range("G72").Value = 2
$ws.Range("H72").Value = 29

$ws.Range("A79").Value = "Kenia"
$ws.Range("B79").Value = 7188
$ws.Range("C79").Value = 247
$ws.Range("D79").Value = 2148
$ws.Range("E79").Value = 4886
$ws.Range("F79").Value = 0
$ws.Range("G79").Value = 2
$ws.Range("H79").Value = 154

$ws.Range("A87").Value = "Guinea"
$ws.Range("B87").Value = 5521
$ws.Range("C87").Value = 71
$ws.Range("D87").Value = 4446
$ws.Range("E87").Value = 1042
$ws.Range("F87").Value = 0
$ws.Range("G87").Value = 0
$ws.Range("H87").Value = 33

$ws.Range("A88").Value = "Gabon"
$ws.Range("B88").Value = 5513
$ws.Range("C88").Value = 0
$ws.Range("D88").Value = 2508
$ws.Range("E88").Value = 2963
$ws.Range("F88").Value = 0
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 42

$ws.Range("A92").Value = "Mauritania"
$ws.Range("B92").Value = 4705
$ws.Range("C92").Value = 99
$ws.Range("D92").Value = 1765
$ws.Range("E92").Value = 2811
$ws.Range("F92").Value = 0
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 129

$ws.Range("A95").Value = "Costa Rica"
$ws.Range("B95").Value = 4311
$ws.Range("C95").Value = 288
$ws.Range("D95").Value = 1657
$ws.Range("E95").Value = 2636
$ws.Range("F95").Value = 0
$ws.Range("G95").Value = 0
$ws.Range("H95").Value = 18

$ws.Range("A96").Value = "Hungria"
$ws.Range("B96").Value = 4172
$ws.Range("C96").Value = 6
$ws.Range("D96").Value = 2752
$ws.Range("E96").Value = 832
$ws.Range("F96").Value = 0
$ws.Range("G96").Value = 1
$ws.Range("H96").Value = 588

$ws.Range("A99").Value = "Estado de Palestina"
$ws.Range("B99").Value = 3334
$ws.Range("C99").Value = 254
$ws.Range("D99").Value = 463
$ws.Range("E99").Value = 2860
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 2
$ws.Range("H99").Value = 11

$ws.Range("A110").Value = "Paraguay"
$ws.Range("B110").Value = 2349
$ws.Range("C110").Value = 46
$ws.Range("D110").Value = 1113
$ws.Range("E110").Value = 1217
$ws.Range("F110").Value = 0
$ws.Range("G110").Value = 0
$ws.Range("H110").Value = 19

$ws.Range("A118").Value = "Guinea-Bisau"
$ws.Range("B118").Value = 1765
$ws.Range("C118").Value = 111
$ws.Range("D118").Value = 676
$ws.Range("E118").Value = 1064
$ws.Range("F118").Value = 0
$ws.Range("G118").Value = 1
$ws.Range("H118").Value = 25

$ws.Range("A119").Value = "Eslovaquia"
$ws.Range("B119").Value = 1720
$ws.Range("C119").Value = 20
$ws.Range("D119").Value = 1466
$ws.Range("E119").Value = 226
$ws.Range("F119").Value = 0
$ws.Range("G119").Value = 0
$ws.Range("H119").Value = 28

$ws.Range("A151").Value = "Togo"
$ws.Range("B151").Value = 671
$ws.Range("C151").Value = 4
$ws.Range("D151").Value = 424
$ws.Range("E151").Value = 233
$ws.Range("F151").Value = 0
$ws.Range("G151").Value = 0
$ws.Range("H151").Value = 14

$ws.Range("A205").Value = "Fiyi"
$ws.Range("B205").Value = 18
$ws.Range("C205").Value = 0
$ws.Range("D205").Value = 18
$ws.Range("E205").Value = 0
$ws.Range("F205").Value = 0
$ws.Range("G205").Value = 0
$ws.Range("H205").Value = 0

$ws.Range("A206").Value = "Dominica"
$ws.Range("B206").Value = 18
$ws.Range("C206").Value = 0
$ws.Range("D206").Value = 18
$ws.Range("E206").Value = 0
$ws.Range("F206").Value = 0
$ws.Range("G206").Value = 0
$ws.Range("H206").Value = 0
